{"js": "// Update the dated math worksheet: the header date moves forward one day,\n// and every two-digit-by-one-digit division problem in the 5x5 table is\n// replaced with that day's new set of problems.\n//\n// Every \"old\" string below occurs exactly once in the starting document,\n// and the pairs are listed in document order, so replacing them in this\n// order can never let one replacement's *new* text be re-matched by a\n// *later* replacement's search (even though a couple of new/old values\n// happen to collide, e.g. \"61\u00f73=\" is both removed early on and produced\n// again later from a different cell).\nconst replacements = [\n  [\"2025-11-10 Monday\", \"2025-11-11 Tuesday\"],\n  [\"35\u00f75=\", \"45\u00f79=\"],\n  [\"16\u00f76=\", \"46\u00f74=\"],\n  [\"20\u00f74=\", \"37\u00f76=\"],\n  [\"44\u00f72=\", \"93\u00f76=\"],\n  [\"13\u00f72=\", \"26\u00f75=\"],\n  [\"61\u00f73=\", \"29\u00f72=\"],\n  [\"66\u00f79=\", \"51\u00f73=\"],\n  [\"79\u00f79=\", \"95\u00f78=\"],\n  [\"38\u00f72=\", \"91\u00f74=\"],\n  [\"17\u00f76=\", \"21\u00f72=\"],\n  [\"90\u00f79=\", \"66\u00f77=\"],\n  [\"47\u00f78=\", \"59\u00f72=\"],\n  [\"94\u00f74=\", \"87\u00f74=\"],\n  [\"60\u00f76=\", \"31\u00f77=\"],\n  [\"61\u00f79=\", \"90\u00f78=\"],\n  [\"66\u00f78=\", \"67\u00f77=\"],\n  [\"40\u00f77=\", \"48\u00f77=\"],\n  [\"64\u00f78=\", \"88\u00f73=\"],\n  [\"49\u00f75=\", \"66\u00f73=\"],\n  [\"31\u00f74=\", \"71\u00f77=\"],\n  [\"33\u00f77=\", \"70\u00f77=\"],\n  [\"76\u00f72=\", \"14\u00f72=\"],\n  [\"50\u00f74=\", \"61\u00f73=\"],\n  [\"15\u00f77=\", \"15\u00f74=\"],\n  [\"61\u00f76=\", \"87\u00f79=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the dated math worksheet: the header date moves forward one day,\n# and every two-digit-by-one-digit division problem in the 5x5 table is\n# replaced with that day's new set of problems.\n#\n# Every \"old\" string below occurs exactly once in the starting document,\n# and the pairs are listed in document order, so replacing them in this\n# order can never let one replacement's *new* text be re-matched by a\n# *later* replacement's search (even though a couple of new/old values\n# happen to collide, e.g. \"61\u00f73=\" is both removed early on and produced\n# again later from a different cell).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-10 Monday\", \"2025-11-11 Tuesday\"),\n    @(\"35\u00f75=\", \"45\u00f79=\"),\n    @(\"16\u00f76=\", \"46\u00f74=\"),\n    @(\"20\u00f74=\", \"37\u00f76=\"),\n    @(\"44\u00f72=\", \"93\u00f76=\"),\n    @(\"13\u00f72=\", \"26\u00f75=\"),\n    @(\"61\u00f73=\", \"29\u00f72=\"),\n    @(\"66\u00f79=\", \"51\u00f73=\"),\n    @(\"79\u00f79=\", \"95\u00f78=\"),\n    @(\"38\u00f72=\", \"91\u00f74=\"),\n    @(\"17\u00f76=\", \"21\u00f72=\"),\n    @(\"90\u00f79=\", \"66\u00f77=\"),\n    @(\"47\u00f78=\", \"59\u00f72=\"),\n    @(\"94\u00f74=\", \"87\u00f74=\"),\n    @(\"60\u00f76=\", \"31\u00f77=\"),\n    @(\"61\u00f79=\", \"90\u00f78=\"),\n    @(\"66\u00f78=\", \"67\u00f77=\"),\n    @(\"40\u00f77=\", \"48\u00f77=\"),\n    @(\"64\u00f78=\", \"88\u00f73=\"),\n    @(\"49\u00f75=\", \"66\u00f73=\"),\n    @(\"31\u00f74=\", \"71\u00f77=\"),\n    @(\"33\u00f77=\", \"70\u00f77=\"),\n    @(\"76\u00f72=\", \"14\u00f72=\"),\n    @(\"50\u00f74=\", \"61\u00f73=\"),\n    @(\"15\u00f77=\", \"15\u00f74=\"),\n    @(\"61\u00f76=\", \"87\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace) -- Wrap:=1 is wdFindContinue, Replace:=2 is\n    # wdReplaceAll (there is only ever one match per query here, since\n    # each old string is unique in the document at the time it is\n    # searched for).\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
